# Linear Regression from Scratch
#
# 1) Delete the duplicate "Training" slide (title "Training" + a
#    picture of the training-loss plot; the deck has two copies of
#    this slide back to back - drop the redundant one).
# 2) Refresh the cached "last saved" datetimeFigureOut date fields
#    (Date Placeholder shapes) on the slide master, the slide layouts
#    that carry one, and the notes master from 9/27/2019 -> 10/3/2019.

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "9/27/2019") {
                $shp.TextFrame.TextRange.Text = "10/3/2019"
            }
        }
    }
}

# Slide master date placeholder.
Update-DatePlaceholders $p.SlideMaster.Shapes

# Every slide layout that has its own date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

# Notes master date placeholder.
Update-DatePlaceholders $p.NotesMaster.Shapes

# Remove the duplicate "Training" slide (slide id 339 / index 15 in the
# original deck). The deck actually has two near-identical "Training"
# slides (same shapes/picture, different internal creationId) - the
# later, redundant one is the one being dropped here, so match on its
# stable SlideID rather than on index/text (which can't tell the two
# apart).
$deleted = $false
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    if ($slide.SlideID -eq 339) {
        $slide.Delete()
        $deleted = $true
        break
    }
}

# Fallback, in case SlideID 339 isn't present for some reason: drop the
# second "Training"-titled slide (the first one is kept).
if (-not $deleted) {
    $seen = 0
    for ($si = 1; $si -le $p.Slides.Count; $si++) {
        $slide = $p.Slides.Item($si)
        $title = ""
        if ($slide.Shapes.Count -ge 1 -and $slide.Shapes.Item(1).HasTextFrame) {
            $title = $slide.Shapes.Item(1).TextFrame.TextRange.Text
        }
        if ($title -eq "Training") {
            $seen = $seen + 1
            if ($seen -eq 2) {
                $slide.Delete()
                break
            }
        }
    }
}
